$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number need to be
# explicitly formatted as Text first, otherwise Excel auto-converts the
# literal (e.g. "322.65") into a floating point number and mangles
# trailing zeros / precision (matches native Excel typed-entry behaviour).
$textPriceCells = @("D5","D6","D8","D10","D11","D13","D14","D17","D19","D20","D22","D23","D24","D26","D27","D28","D30","D32","D33","D34","D37","D38","D40","D42","D43","D44","D49","D50","D51")
foreach ($c in $textPriceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "46.838.69"
$ws.Range("D3").Value = "2.481.13"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "322.65"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "104.20"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "36.78"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "18.20"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "7.15"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "2.870.71"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "2.523.25"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("D17").Value = "0.839"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "46.797.42"
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "6.55"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").Value = "0.0₃0930"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "70.49"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "249.84"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").Value = "2.34"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "26.02"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").Value = "34.99"
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").Value = "49.36"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "19.50"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").Value = "5.29"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "4.55"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "122.66"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "21.31"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "0.0294"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "1.950.72"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "9.11"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").Value = "  +13.50%  "
$ws.Range("D51").Value = "78.42"
$ws.Range("E51").Value = "  +2.81%  "
